$wb = $excel.ActiveWorkbook

# The "TestCoverageMatrix" sheet had a leftover "hello world" value (with a
# highlight fill) sitting in A1. Clear it out completely - value and
# formatting - so the header row starts cleanly at B1, matching the real
# test-id headers.
$ws1 = $wb.Worksheets.Item("TestCoverageMatrix")
$ws1.Range("A1").Clear()

# Move the active selection on that sheet, mirroring where the author's
# cursor ended up after the edit.
$ws1.Range("G13").Select()

$wb.Save()
